# BehaviorScenario.xlsx edit: add a new household type (5) and the
# corresponding wfh_share scenarios (0, 0.2, 0.4, 0.6, 0.8, 1), and
# refresh the wfh_share step for every existing household type from
# {0, 0.5, 1} to {0, 0.2, 0.4, 0.6, 0.8, 1}.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the built-in "Standard" cell style to its English name "Normal"
# (workbook was re-saved from a non-English Excel build originally).
$normalStyle = $wb.Styles.Item(1)
if ($normalStyle.Name -eq "Standard") {
    $normalStyle.Name = "Normal"
}

$wfhSteps = @(0, 0.2, 0.4, 0.6, 0.8, 1)

$scenario = 1
$row = 2
for ($householdType = 1; $householdType -le 5; $householdType++) {
    foreach ($wfh in $wfhSteps) {
        $ws.Cells.Item($row, 1).Value = $scenario
        $ws.Cells.Item($row, 2).Value = $householdType
        $ws.Cells.Item($row, 3).Value = $wfh
        $scenario++
        $row++
    }
}

# Move the selection the way the saved file shows (no frozen/scrolled
# top-left cell, active cell sits on the first row of the last block).
[void]$ws.Range("A1").Select()
[void]$ws.Range("E13").Select()
